$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" (col E) / "Valor Mora" (col F) data block occupies
# rows 16-35. The previous account-statement periods are being removed and
# replaced with the new ones: the block is reversed top-to-bottom (row 16
# now shows what used to be on row 35, etc.), so read the current values
# first and write them back in reverse order.

$firstRow = 16
$lastRow = 35

$eVals = @()
$fVals = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $eVals += $ws.Cells.Item($r, 5).Value()
    $fVals += $ws.Cells.Item($r, 6).Value()
}

$eVals = $eVals[($eVals.Count - 1)..0]
$fVals = $fVals[($fVals.Count - 1)..0]

$i = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = $eVals[$i]
    $ws.Cells.Item($r, 6).Value = $fVals[$i]
    $i++
}
